$d = $word.ActiveDocument

# --- Change 1: merge the split runs that make up the "Journal Entry 4" heading
# into a single run, matching how Word coalesces a retyped heading.
$d.Content.Find.Execute("## Journal Entry 4, Mod 6 " + [char]0x2013 + " 2/10/2025", $true, $false, $false, $false, $false, $true, 1, $false, "## Journal Entry 4, Mod 6 " + [char]0x2013 + " 2/10/2025", 2) | Out-Null

# --- Change 2: insert the new "Journal Entry 5" heading paragraph and body
# paragraph right before the trailing empty paragraph at the end of the document.
$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">## Journal Entry </w:t></w:r><w:r><w:t>5</w:t></w:r><w:r><w:t xml:space="preserve">, Mod </w:t></w:r><w:r><w:t>7</w:t></w:r><w:r><w:t xml:space="preserve"> – 2/1</w:t></w:r><w:r><w:t>2</w:t></w:r><w:r><w:t>/2025</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t xml:space="preserve">Chapter 15 focused on Git GUI.  This implementation of git is closer to how I first used git during my SDEV and Computer Science courses.  It was gratifying to see how much I have learned and how much I understand about what is taking place “under the hood” so to speak, when working with more advanced software development tools.  I appreciate the ability that I have gained to take the git commands that I use in Git GUI and break them down in terms of several command lines in git bash.  Tools like Git GUI will now be much less mysterious to me and understanding the command lines will help me resolve issues that might otherwise be difficult to address.  As far as chapter 16 goes, I feel there </w:t></w:r><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">was a great deal of overlap with what I had already done in module six, because I had already used git stash to resolve conflicts.  Git stash will continue to be one of the most important commands that I use moving forward, as I will try my best to follow the author’s </w:t></w:r><w:r><w:t>admonition</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">to avoid reset when possible.  When resetting becomes necessary, I will make sure to review the difference between hard and soft resets.  Another command which I feel will be important in chapter 16 is the checkout feature.  Although I have seen it before, its use in chapter 16 is related to </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>to</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> files and not commits and the results are somewhat different in this context.  </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$n = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($n)
$insertRange = $lastPara.Range
$insertRange.InsertParagraphBefore()

$targetPara = $d.Paragraphs($n)
$targetRange = $targetPara.Range
$targetRange.InsertXML($xml)
